$wb = $excel.ActiveWorkbook

# Sheet "block_0" (first sheet) - row 2 values for C:F
$ws1 = $wb.Worksheets.Item("block_0")
$ws1.Range("C2").Value = [double]"3.315138836796488e-08"
$ws1.Range("D2").Value = [double]"0.001253061629053576"
$ws1.Range("E2").Value = [double]"1.726506744514911e-08"
$ws1.Range("F2").Value = [double]"2.260962129280777e-05"

# Sheet "block_1" (second sheet) - row 2 values for C:F
$ws2 = $wb.Worksheets.Item("block_1")
$ws2.Range("C2").Value = [double]"2.311383758869801e-09"
$ws2.Range("D2").Value = [double]"0.005861695524878178"
$ws2.Range("E2").Value = [double]"1.329795547982959e-09"
$ws2.Range("F2").Value = [double]"0.0001357465553551715"
